# Update F-column (view/heat count) figures across the four sheets of the
# 广州-漫展信息 workbook, as published to gh-pages at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 21409
$ws1.Range("F3").Value = 3302
$ws1.Range("F5").Value = 622
$ws1.Range("F11").Value = 139
$ws1.Range("F12").Value = 569
$ws1.Range("F14").Value = 349
$ws1.Range("F15").Value = 37
$ws1.Range("F16").Value = 453
$ws1.Range("F17").Value = 187
$ws1.Range("F19").Value = 31
$ws1.Range("F20").Value = 80
$ws1.Range("F21").Value = 151

# 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 169

# 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6171
$ws3.Range("F3").Value = 727
$ws3.Range("F4").Value = 728
$ws3.Range("F5").Value = 1730
$ws3.Range("F6").Value = 84

# 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6171
$ws4.Range("F3").Value = 727
$ws4.Range("F4").Value = 728
$ws4.Range("F5").Value = 1730
$ws4.Range("F6").Value = 21409
$ws4.Range("F7").Value = 3302
$ws4.Range("F10").Value = 84
$ws4.Range("F11").Value = 622
$ws4.Range("F20").Value = 139
$ws4.Range("F23").Value = 569
$ws4.Range("F27").Value = 349
$ws4.Range("F28").Value = 169
$ws4.Range("F29").Value = 37
$ws4.Range("F30").Value = 453
$ws4.Range("F32").Value = 187
$ws4.Range("F36").Value = 31
$ws4.Range("F37").Value = 80
$ws4.Range("F43").Value = 151

$wb.Save()
